# Update the dsafnet metrics worksheet with the new training-run data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ---
# Insert the new "learning_rate" column before "epoch_time" and shift epoch_time to H1.
$ws.Range("A1").Value = "epoch"
$ws.Range("B1").Value = "train_loss"
$ws.Range("C1").Value = "val_loss"
$ws.Range("D1").Value = "val_mae"
$ws.Range("E1").Value = "val_r2"
$ws.Range("F1").Value = "val_rmse"
$ws.Range("G1").Value = "learning_rate"
$ws.Range("H1").Value = "epoch_time"

# Copy the header style (bold/centered/bordered) from G1 (already styled) to the new H1 cell.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows ---
$data = @(
    @(1, 0.9033397331730119, 0.8458161743650658, 6.391531875526314, 0.1672741257265296, 8.945469716399325, 0.0005, 71.36291837692261),
    @(2, 0.8928861508789859, 0.900063429837065, 7.298406696147556, 0.05504867505692453, 9.523915326063696, 0.0005, 58.57611536979675),
    @(3, 0.8904605047332161, 0.8825040088079323, 6.971020179235587, 0.09011494626397053, 9.345961646774244, 0.0005, 52.9040858745575),
    @(4, 0.8889403258841442, 0.8606795145820689, 6.685294652557366, 0.1366448900075727, 9.107482272116034, 0.0005, 53.81499409675598),
    @(5, 0.8884218232122763, 0.871254037527264, 6.839464001744581, 0.1149594803934071, 9.218530078513517, 0.0005, 53.7888720035553),
    @(6, 0.8870917813130709, 0.8626975261211941, 6.680427262495004, 0.1324329599492154, 9.1288671856053, 0.00025, 57.75950288772583),
    @(7, 0.8867508136341025, 0.8606720745404399, 6.688650655584258, 0.1361438733672364, 9.110516361816915, 0.00025, 61.31318712234497)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($value in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $value
        $colIndex = $colIndex + 1
    }
    $rowIndex = $rowIndex + 1
}
